# added the extra ov times
# Update the OV (public transport) travel-time matrix on Sheet1 with the
# newly measured values, then leave the selection over the data block
# (F2:Y21) as the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35
$ws.Range("J2").Value = 30
$ws.Range("O2").Value = 46
$ws.Range("Q2").Value = 31
$ws.Range("T2").Value = 45
$ws.Range("U2").Value = 43

$ws.Range("F3").Value = 35
$ws.Range("I3").Value = 42
$ws.Range("T3").Value = 40
$ws.Range("U3").Value = 30

$ws.Range("J4").Value = 41
$ws.Range("M4").Value = 22

$ws.Range("G5").Value = 42
$ws.Range("U5").Value = 40
$ws.Range("X5").Value = 17

$ws.Range("F6").Value = 30
$ws.Range("H6").Value = 41
$ws.Range("K6").Value = 22
$ws.Range("L6").Value = 27
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 16
$ws.Range("O6").Value = 22
$ws.Range("P6").Value = 8
$ws.Range("Q6").Value = 30
$ws.Range("U6").Value = 38
$ws.Range("X6").Value = 17
$ws.Range("Y6").Value = 26

$ws.Range("J7").Value = 22
$ws.Range("O7").Value = 26
$ws.Range("Q7").Value = 24

$ws.Range("J8").Value = 27
$ws.Range("M8").Value = 15
$ws.Range("R8").Value = 27
$ws.Range("V8").Value = 14
$ws.Range("W8").Value = 25

$ws.Range("H9").Value = 22
$ws.Range("J9").Value = 30
$ws.Range("L9").Value = 15
$ws.Range("R9").Value = 33
$ws.Range("V9").Value = 21
$ws.Range("W9").Value = 11

$ws.Range("J10").Value = 16
$ws.Range("X10").Value = 20

$ws.Range("F11").Value = 46
$ws.Range("J11").Value = 22
$ws.Range("K11").Value = 26

$ws.Range("J12").Value = 8
$ws.Range("X12").Value = 18
$ws.Range("Y12").Value = 20

$ws.Range("F13").Value = 31
$ws.Range("J13").Value = 30
$ws.Range("K13").Value = 24
$ws.Range("R13").Value = 25
$ws.Range("S13").Value = 22
$ws.Range("T13").Value = 36

$ws.Range("L14").Value = 27
$ws.Range("M14").Value = 33
$ws.Range("Q14").Value = 25
$ws.Range("V14").Value = 27
$ws.Range("W14").Value = 38

$ws.Range("Q15").Value = 22

$ws.Range("F16").Value = 45
$ws.Range("G16").Value = 40
$ws.Range("Q16").Value = 36

$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 30
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 38
$ws.Range("Y17").Value = 50

$ws.Range("L18").Value = 14
$ws.Range("M18").Value = 21
$ws.Range("R18").Value = 27
$ws.Range("W18").Value = 25

$ws.Range("L19").Value = 25
$ws.Range("M19").Value = 11
$ws.Range("R19").Value = 38
$ws.Range("V19").Value = 25
$ws.Range("X19").Value = 39

$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 17
$ws.Range("N20").Value = 20
$ws.Range("P20").Value = 18
$ws.Range("W20").Value = 39
$ws.Range("Y20").Value = 32

$ws.Range("J21").Value = 26
$ws.Range("P21").Value = 20
$ws.Range("U21").Value = 50
$ws.Range("X21").Value = 32

# Match the author's final selection in the saved workbook.
$ws.Range("F2:Y21").Select()
